$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.344.04'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.686.71'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.19%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.90%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.44'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5459'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.74%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.011'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.87%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2722'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.88%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06442'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.01'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07660'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.79%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.702.20'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.532'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5805'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008333'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.09'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.384.53'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.941'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.011'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.80%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.96'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.33'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.29%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.222'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.012'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.92%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.86'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1315'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +5.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.884'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.70'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.90%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06342'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.411'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.329'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.73%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.578'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.573'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.673'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.040'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6159'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.81%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.411'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.718'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.55%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.230'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.114.50'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01627'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8793'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.016'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.19'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.00%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.839.77'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.24'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.191'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.009'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.48%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05271'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4306'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.028'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.19%  '

